$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# --- Step 1: Insert a new row at 70, shifting rows 70-129 down to 71-130 ---
$ws.Rows("70:70").Insert()

# --- Step 2: Fix up B70 style to match the other field rows (copy format from B71) ---
$ws.Range("B71").Copy()
$ws.Range("B70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 3: Set the new row 70 label ---
$ws.Range('A70').Value = 'VSTAT License File'

# --- Step 4: Shift existing comments (rows 70-129) down to (71-130). ---
# Process from the bottom row upward so we never clobber a comment we still need to read.
$c = $ws.Range('A129').Comment
$c.Delete()
$newC = $ws.Range('A130').AddComment('Allowing VSD in-place upgrade during Installation [default: False]')
$c = $ws.Range('A128').Comment
$c.Delete()
$newC = $ws.Range('A129').AddComment('Events for which alerts should not be sent. One string can be used to hold multiple events, separated by commas')
$c = $ws.Range('A127').Comment
$c.Delete()
$newC = $ws.Range('A128').AddComment('Specific events for which alerts should be sent. One string can be used to hold multiple events, separated by commas')
$c = $ws.Range('A126').Comment
$c.Delete()
$newC = $ws.Range('A127').AddComment('Destination email address for monit alerts')
$c = $ws.Range('A125').Comment
$c.Delete()
$newC = $ws.Range('A126').AddComment('Email message for alert emails. Overrides monit default alert message')
$c = $ws.Range('A124').Comment
$c.Delete()
$newC = $ws.Range('A125').AddComment('Email subject for alert emails. Overrides monit default alert subject')
$c = $ws.Range('A123').Comment
$c.Delete()
$newC = $ws.Range('A124').AddComment('Email address to reply to monit alert emails')
$c = $ws.Range('A122').Comment
$c.Delete()
$newC = $ws.Range('A123').AddComment('Email address from which monit alerts will be sent')
$c = $ws.Range('A121').Comment
$c.Delete()
$newC = $ws.Range('A122').AddComment('Enables use of monit eventqueue to store alerts if email alerts fail to send [default: True]')
$c = $ws.Range('A120').Comment
$c.Delete()
$newC = $ws.Range('A121').AddComment('Encryption to be used when sending monit alerts via email')
$c = $ws.Range('A119').Comment
$c.Delete()
$newC = $ws.Range('A120').AddComment('Port on mail server to be used for monit alerts [default: 25]')
$c = $ws.Range('A118').Comment
$c.Delete()
$newC = $ws.Range('A119').AddComment('Address of the mail server to be used to receive monit alerts via email')
$c = $ws.Range('A116').Comment
$c.Delete()
$newC = $ws.Range('A117').AddComment('List of destination email addresses (List items separated by comma.)')
$c = $ws.Range('A115').Comment
$c.Delete()
$newC = $ws.Range('A116').AddComment('Email address from which health report will be sent')
$c = $ws.Range('A114').Comment
$c.Delete()
$newC = $ws.Range('A115').AddComment('Port to be used on the SMTP Server [default: 25]')
$c = $ws.Range('A113').Comment
$c.Delete()
$newC = $ws.Range('A114').AddComment('Address of SMTP server to be used if emailing health results')
$c = $ws.Range('A111').Comment
$c.Delete()
$newC = $ws.Range('A112').AddComment('Skip tasks and playbooks (List items separated by comma.)')
$c = $ws.Range('A110').Comment
$c.Delete()
$newC = $ws.Range('A111').AddComment('List of hooks files (List items separated by comma.)')
$c = $ws.Range('A108').Comment
$c.Delete()
$newC = $ws.Range('A109').AddComment('Enterprise name used for authentication with VCIN. Required for tasks like VRS-E upgrade (through VCIN) [default: csp]')
$c = $ws.Range('A107').Comment
$c.Delete()
$newC = $ws.Range('A108').AddComment('VCIN URL used for API interaction. Required for tasks like VRS-E upgrade (through VCIN) [default: https://(vcin_ip_address):8443]')
$c = $ws.Range('A106').Comment
$c.Delete()
$newC = $ws.Range('A107').AddComment('Enterprise name used for authentication with VSD Architect. Required for tasks during Upgrade, Health Checks etc [default: csp]')
$c = $ws.Range('A105').Comment
$c.Delete()
$newC = $ws.Range('A106').AddComment('VSD Architect URL. Required for tasks during Upgrade, Health Checks etc [default: https://(vsd_fqdn):8443]')
$c = $ws.Range('A103').Comment
$c.Delete()
$newC = $ws.Range('A104').AddComment('Number of CPU''s for Webfilter vm. Valid only for KVM deployments [default: 2]')
$c = $ws.Range('A102').Comment
$c.Delete()
$newC = $ws.Range('A103').AddComment('Number of CPU''s for Portal vm. Valid only for KVM deployments [default: 6]')
$c = $ws.Range('A101').Comment
$c.Delete()
$newC = $ws.Range('A102').AddComment('Number of CPU''s for VCIN. Valid only for KVM deployments [default: 6]')
$c = $ws.Range('A100').Comment
$c.Delete()
$newC = $ws.Range('A101').AddComment('Number of CPU''s for NUH. Valid only for KVM deployments [default: 2]')
$c = $ws.Range('A98').Comment
$c.Delete()
$newC = $ws.Range('A99').AddComment('Valid for only KVM and VCenter deployments. Number of CPU''s for VNSUTIL. [default: 2]')
$c = $ws.Range('A97').Comment
$c.Delete()
$newC = $ws.Range('A98').AddComment('Valid for only KVM and VCenter deployments. Number of CPU''s for VSTAT. [default: 6]')
$c = $ws.Range('A96').Comment
$c.Delete()
$newC = $ws.Range('A97').AddComment('Valid for only KVM and VCenter deployments. Number of CPU''s for VSC. [default: 6]')
$c = $ws.Range('A95').Comment
$c.Delete()
$newC = $ws.Range('A96').AddComment('Valid for only KVM and VCenter deployments. Number of CPU''s for VSD. [default: 6]')
$c = $ws.Range('A93').Comment
$c.Delete()
$newC = $ws.Range('A94').AddComment('Amount of Portal RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]')
$c = $ws.Range('A92').Comment
$c.Delete()
$newC = $ws.Range('A93').AddComment('Amount of Webfilter RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 8]')
$c = $ws.Range('A91').Comment
$c.Delete()
$newC = $ws.Range('A92').AddComment('Amount of NUH RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 8]')
$c = $ws.Range('A90').Comment
$c.Delete()
$newC = $ws.Range('A91').AddComment('Amount of VCIN RAM to allocate, in GB. Valid only for KVM deployments. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]')
$c = $ws.Range('A88').Comment
$c.Delete()
$newC = $ws.Range('A89').AddComment('Valid for only KVM and VCenter deployments. Amount of VSTAT RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 16]')
$c = $ws.Range('A87').Comment
$c.Delete()
$newC = $ws.Range('A88').AddComment('Valid for only KVM and VCenter deployments. Amount of VSC RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 4]')
$c = $ws.Range('A86').Comment
$c.Delete()
$newC = $ws.Range('A87').AddComment('Valid for only KVM and VCenter deployments. Amount of VSD RAM to allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 24]')
$c = $ws.Range('A84').Comment
$c.Delete()
$newC = $ws.Range('A85').AddComment('Amount of NSGV disk space to pre-allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments should not modify this value. [default: 4]')
$c = $ws.Range('A83').Comment
$c.Delete()
$newC = $ws.Range('A84').AddComment('Amount of VCIN disk space to pre-allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments should not modify this value. [default: 285]')
$c = $ws.Range('A82').Comment
$c.Delete()
$newC = $ws.Range('A83').AddComment('Amount of Portal disk space to pre-allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments should not modify this value. [default: 16]')
$c = $ws.Range('A81').Comment
$c.Delete()
$newC = $ws.Range('A82').AddComment('Amount of VSTAT disk space to pre-allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments should not modify this value. [default: 100]')
$c = $ws.Range('A80').Comment
$c.Delete()
$newC = $ws.Range('A81').AddComment('Amount of VSC disk space to pre-allocate, in GB. The only valid values are 0 and 1. When undefined or 0, file size allocation will be skipped. Production deployments should set this value to 1. [default: 0]')
$c = $ws.Range('A79').Comment
$c.Delete()
$newC = $ws.Range('A80').AddComment('Amount of VSD disk space to pre-allocate, in GB. Note: Values smaller than the default are for lab and PoC only. Production deployments must use a value greater than or equal to the default. [default: 285]')
$c = $ws.Range('A77').Comment
$c.Delete()
$newC = $ws.Range('A78').AddComment('IP Address or Hostname of the SSH host if one is used [default: sshhost]')
$c = $ws.Range('A76').Comment
$c.Delete()
$newC = $ws.Range('A77').AddComment('Username of the SSH proxy host if one is used [default: root]')
$c = $ws.Range('A75').Comment
$c.Delete()
$newC = $ws.Range('A76').AddComment('Proxy URL to be used if Yum repositories cannot be directly reached [default: NONE]')
$c = $ws.Range('A74').Comment
$c.Delete()
$newC = $ws.Range('A75').AddComment('Flag to indicate whether to perform a Yum update on VSTAT during the installation [default: False]')
$c = $ws.Range('A73').Comment
$c.Delete()
$newC = $ws.Range('A74').AddComment('Flag to indicate whether to perform a Yum update on VSD during the installation [default: True]')
$c = $ws.Range('A71').Comment
$c.Delete()
$newC = $ws.Range('A72').AddComment('Path to the license file for the NUH including the file name [default: ]')
$c = $ws.Range('A70').Comment
$c.Delete()
$newC = $ws.Range('A71').AddComment('Path to the license file for the SD-WAN Portal including the file name [default: ]')

# --- Step 5: Add the brand new VSTAT License File comment at A70 ---
$ws.Range('A70').AddComment('Optional License File for Elasticsearch [default: ]')

Write-Host "Done"
